$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.159.15"
$ws.Range("E2").Value = "  +4.70%  "

$ws.Range("D3").Value = "2.240.25"
$ws.Range("E3").Value = "  +4.44%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.50"
$ws.Range("E5").Value = "  +6.12%  "

$ws.Range("E6").Value = "  +2.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.12"
$ws.Range("E7").Value = "  +8.86%  "

$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  +6.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.09"
$ws.Range("E10").Value = "  +7.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0932"
$ws.Range("E11").Value = "  +4.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.92"
$ws.Range("E12").Value = "  +5.16%  "

$ws.Range("E13").Value = "  +2.08%  "

$ws.Range("D14").Value = "2.577.94"
$ws.Range("E14").Value = "  +4.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.62"
$ws.Range("E15").Value = "  +1.30%  "

$ws.Range("D16").Value = "2.250.84"
$ws.Range("E16").Value = "  +5.31%  "

$ws.Range("E17").Value = "  +2.26%  "

$ws.Range("D18").Value = "43.064.22"
$ws.Range("E18").Value = "  +5.06%  "

$ws.Range("E19").Value = "  +6.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.20"
$ws.Range("E20").Value = "  +3.11%  "

$ws.Range("E21").Value = "  +5.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.20"
$ws.Range("E22").Value = "  +17.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.73"
$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "229.57"
$ws.Range("E24").Value = "  +2.56%  "

$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("E26").Value = "  +2.98%  "

$ws.Range("E27").Value = "  +2.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.32"
$ws.Range("E28").Value = "  +28.71%  "

$ws.Range("E29").Value = "  +5.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.23"
$ws.Range("E30").Value = "  +3.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.83"
$ws.Range("E31").Value = "  +2.18%  "

$ws.Range("E32").Value = "  +3.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0801"
$ws.Range("E33").Value = "  +7.09%  "

$ws.Range("E34").Value = "  +4.85%  "

$ws.Range("E35").Value = "  +2.52%  "

$ws.Range("E36").Value = "  +10.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.47"
$ws.Range("E37").Value = "  +10.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0332"
$ws.Range("E38").Value = "  +19.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.12"
$ws.Range("E39").Value = "  +13.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.11"
$ws.Range("E40").Value = "  +4.35%  "

$ws.Range("E41").Value = "  +11.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.43"
$ws.Range("E42").Value = "  +3.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.43"
$ws.Range("E43").Value = "  +3.85%  "

$ws.Range("B44").Value = "WOONetwork"
$ws.Range("C44").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.487"
$ws.Range("E44").Value = "  +35.16%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.66"
$ws.Range("E45").Value = "  +6.49%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.11"
$ws.Range("E46").Value = "  +8.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0991"
$ws.Range("E47").Value = "  +4.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  +13.47%  "

$ws.Range("E49").Value = "  +4.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  +4.89%  "

$ws.Range("E51").Value = "  +3.84%  "
